$d = $word.ActiveDocument

function New-PkgXml($bodyInnerXml) {
    return @"
<?xml version='1.0'?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"
                  xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing"
                  xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"
                  xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"
                  xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">
        <w:body>
$bodyInnerXml
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@
}

# ------------------------------------------------------------------
# 1) Move the "_GoBack" bookmark from the end of the "Project:" line
#    to the very start of the document (right after the first
#    paragraph's pPr, before its first run). The first paragraph
#    holds an inline picture, so we rebuild it verbatim (picture +
#    all) with the bookmark spliced in right after <w:pPr>, instead
#    of using Bookmarks.Add (which mis-anchors the end marker when a
#    drawing object sits in the same paragraph as the collapsed
#    range).
# ------------------------------------------------------------------
$p1Body = @"
<w:p>
  <w:pPr>
    <w:pStyle w:val="Date"/>
    <w:jc w:val="left"/>
    <w:outlineLvl w:val="0"/>
    <w:rPr><w:bCs/><w:sz w:val="24"/></w:rPr>
  </w:pPr>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
  <w:r>
    <w:rPr><w:bCs/><w:noProof/><w:sz w:val="24"/></w:rPr>
    <w:drawing>
      <wp:inline distT="0" distB="0" distL="0" distR="0">
        <wp:extent cx="3564613" cy="762647"/>
        <wp:effectExtent l="19050" t="0" r="0" b="0"/>
        <wp:docPr id="2" name="Picture 1" descr="BCS Logo 092310.bmp"/>
        <wp:cNvGraphicFramePr>
          <a:graphicFrameLocks noChangeAspect="1"/>
        </wp:cNvGraphicFramePr>
        <a:graphic>
          <a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture">
            <pic:pic>
              <pic:nvPicPr>
                <pic:cNvPr id="0" name="BCS Logo 092310.bmp"/>
                <pic:cNvPicPr/>
              </pic:nvPicPr>
              <pic:blipFill>
                <a:blip r:embed="rId7"/>
                <a:stretch><a:fillRect/></a:stretch>
              </pic:blipFill>
              <pic:spPr>
                <a:xfrm>
                  <a:off x="0" y="0"/>
                  <a:ext cx="3564613" cy="762647"/>
                </a:xfrm>
                <a:prstGeom prst="rect"><a:avLst/></a:prstGeom>
              </pic:spPr>
            </pic:pic>
          </a:graphicData>
        </a:graphic>
      </wp:inline>
    </w:drawing>
  </w:r>
</w:p>
"@
$d.Paragraphs.Item(1).Range.InsertXML((New-PkgXml $p1Body))

# ------------------------------------------------------------------
# 2) "Date:  `DATE~" -> "Date:  09/23/2015"
#    (keep the existing "Date:  " run untouched; replace the three
#    runs that made up the placeholder with a single plain run)
# ------------------------------------------------------------------
$p4Body = @"
<w:p>
  <w:pPr><w:pStyle w:val="Date"/></w:pPr>
  <w:r><w:t xml:space="preserve">Date:  </w:t></w:r>
  <w:r><w:t>09/23/2015</w:t></w:r>
</w:p>
"@
$d.Paragraphs.Item(4).Range.InsertXML((New-PkgXml $p4Body))

# ------------------------------------------------------------------
# 3) "Project:	 `SHORT~" (+ bookmark) -> "Project:	 Short Title"
#    (keep the field-code structure around "Project:" intact, drop
#    the old `_GoBack` bookmark that used to sit here)
# ------------------------------------------------------------------
$p5Body = @"
<w:p>
  <w:pPr><w:pStyle w:val="Date"/></w:pPr>
  <w:r><w:fldChar w:fldCharType="begin"/></w:r>
  <w:r><w:instrText xml:space="preserve"> AUTOTEXTLIST  </w:instrText></w:r>
  <w:r><w:fldChar w:fldCharType="separate"/></w:r>
  <w:r><w:t>Project:</w:t></w:r>
  <w:r><w:tab/></w:r>
  <w:r><w:fldChar w:fldCharType="end"/></w:r>
  <w:r><w:t xml:space="preserve"> </w:t></w:r>
  <w:r><w:t>Short Title</w:t></w:r>
</w:p>
"@
$d.Paragraphs.Item(5).Range.InsertXML((New-PkgXml $p5Body))

# ------------------------------------------------------------------
# 4) BCS paragraph: wrap the word "to" in proofErr gramStart/gramEnd
#    markers (splitting the single run into three runs).
# ------------------------------------------------------------------
$p7Body = @"
<w:p>
  <w:pPr>
    <w:pStyle w:val="BodyText"/>
    <w:spacing w:line="360" w:lineRule="auto"/>
    <w:rPr><w:rFonts w:cs="Arial"/></w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr><w:rFonts w:cs="Arial"/></w:rPr>
    <w:t xml:space="preserve">BCS service offers installation of commercial and industrial grade cabling solutions.  BCS installations include inside plant cabling of Category 5e/6, fiber optics, coaxial, security and control cabling.  BCS also provides outside plant cabling for high pair copper, fiber optics, and coaxial. All projects are supported by a Registered Certified Distribution Designer (BICSI, RCDD) and Certified BICSI Installers </w:t>
  </w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr><w:rFonts w:cs="Arial"/></w:rPr>
    <w:t>to</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:rPr><w:rFonts w:cs="Arial"/></w:rPr>
    <w:t xml:space="preserve"> ensure strict adherence to industry standards and a successful design and construction team. </w:t>
  </w:r>
</w:p>
"@
$d.Paragraphs.Item(7).Range.InsertXML((New-PkgXml $p7Body))

Write-Host "Done."
